# Generate Report for Handoff
#
# Updates the "Latest Handoff Datetime" (column D) on the per-locale
# status sheets (zh-cn, de-de) for every row whose status is
# "In Translation" (rows 4, 6, 7, 8, 9, 10) to reflect a freshly
# generated handoff report timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(4, 6, 7, 8, 9, 10)

foreach ($r in $rows) {
    $zhcn.Range("D$r").Value = "2016-03-07 03:01:13"
}

foreach ($r in $rows) {
    $dede.Range("D$r").Value = "2016-03-07 03:01:27"
}
